$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feb")

# Row 3: RAMESHVAR FARKANDE - Link Triggered count changes 2 -> 4
$ws.Range("B3").Value = 4

# Row 4: ANIKET DESHMUKH - Link Triggered count changes 45 -> 46
$ws.Range("B4").Value = 46

# Row 6: ASHUTOSH GURAV - Link Triggered 29 -> 31, CC/1000 34.48 -> 32.26
$ws.Range("B6").Value = 31
$ws.Range("E6").Value = 32.26

# Row 7: AVINASH KAMBLE - Link Triggered 39 -> 43, CC/1000 51.28 -> 46.51
$ws.Range("B7").Value = 43
$ws.Range("E7").Value = 46.51

# Row 9: now DATTA SHEJAV (new SA) instead of DNYANESHWAR GAWADE; Response/Concern Count become blank; Division -> AMRAVATI
$ws.Range("A9").Value = "DATTA SHEJAV"
$ws.Range("B9").Value = 2
$ws.Range("C9").ClearContents()
$ws.Range("F9").Value = "AMRAVATI"

# Row 10: now DNYANESHWAR GAWADE (shifted down); Link Triggered 18, Response 2
$ws.Range("A10").Value = "DNYANESHWAR GAWADE"
$ws.Range("B10").Value = 18
$ws.Range("C10").Value = 2

# Row 11: now JAVED RAMPURE (shifted down); Link Triggered 15, Response 1, Division -> WAGHOLI
$ws.Range("A11").Value = "JAVED RAMPURE"
$ws.Range("B11").Value = 15
$ws.Range("C11").Value = 1
$ws.Range("F11").Value = "WAGHOLI"

# Row 12: now MANOJ PATIL (shifted down); Link Triggered 20, Division -> KOLHAPUR_WS
$ws.Range("A12").Value = "MANOJ PATIL"
$ws.Range("B12").Value = 20
$ws.Range("F12").Value = "KOLHAPUR_WS"

# Row 13: now MOHSIN ALI (shifted down); Link Triggered 22
$ws.Range("A13").Value = "MOHSIN ALI"
$ws.Range("B13").Value = 22

# Row 14: now SANJAY RAMKELKAR (shifted down); Link Triggered 10, Response/Concern Count blank, CC/1000 0
$ws.Range("A14").Value = "SANJAY RAMKELKAR"
$ws.Range("B14").Value = 10
$ws.Range("C14").ClearContents()
$ws.Range("D14").ClearContents()
$ws.Range("E14").Value = 0

# Row 15 (new row): UJJWAL MAHAJAN
$ws.Range("A15").Value = "UJJWAL MAHAJAN"
$ws.Range("B15").Value = 45
$ws.Range("C15").Value = 1
$ws.Range("D15").Value = 1
$ws.Range("E15").Value = 22.22
$ws.Range("F15").Value = "NAGPUR_KAMPTHEE ROAD"

# Row 16 (was row 15): VAIBHAV PANCHAL, now shifted to row 16, Division -> AMRAVATI
$ws.Range("A16").Value = "VAIBHAV PANCHAL"
$ws.Range("B16").Value = 42
$ws.Range("C16").Value = 3
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 23.81
$ws.Range("F16").Value = "AMRAVATI"

# Column F width narrows to fit the new shorter content
$ws.Columns.Item(6).ColumnWidth = 13

# Update selection to match the saved state of the workbook after editing
$ws.Range("F8").Select()
